$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number must be forced
# to remain text (matching the sheet's existing text-formatted Price column),
# otherwise Excel's normal type inference would convert them to numbers.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '51.615.80'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.940.27'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '378.55'
$ws.Range('E5').Value = '  +6.74%  '
$ws.Range('D6').Value = '104.48'
$ws.Range('E6').Value = '  -1.77%  '
$ws.Range('D7').Value = '0.542'
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '0.589'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('D10').Value = '37.07'
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '0.0839'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').Value = '18.45'
$ws.Range('E13').Value = '  -2.65%  '
$ws.Range('D14').Value = '3.397.08'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '7.41'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '2.932.44'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '0.946'
$ws.Range('E17').Value = '  -4.87%  '
$ws.Range('D18').Value = '51.550.39'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').Value = '7.35'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').Value = '0.0₃0952'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').Value = '68.46'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').Value = '262.50'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '2.82'
$ws.Range('E25').Value = '  +4.39%  '
$ws.Range('D26').Value = '7.40'
$ws.Range('E26').Value = '  +18.21%  '
$ws.Range('D27').Value = '4.13'
$ws.Range('E27').Value = '  -4.92%  '
$ws.Range('E29').Value = '  -5.04%  '
$ws.Range('D30').Value = '7.37'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D31').Value = '25.88'
$ws.Range('E31').Value = '  -3.11%  '
$ws.Range('E32').Value = '  -5.51%  '
$ws.Range('D33').Value = '9.87'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = '52.42'
$ws.Range('E34').Value = '  +2.98%  '
$ws.Range('E35').Value = '  -3.21%  '
$ws.Range('D36').Value = '34.27'
$ws.Range('E36').Value = '  -3.93%  '
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('E39').Value = '  -6.75%  '
$ws.Range('D40').Value = '17.14'
$ws.Range('E40').Value = '  -1.31%  '
$ws.Range('E41').Value = '  -6.39%  '
$ws.Range('E42').Value = '  -5.00%  '
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('D44').Value = '124.03'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '21.79'
$ws.Range('E45').Value = '  -4.98%  '
$ws.Range('D46').Value = '2.07'
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('D47').Value = '0.277'
$ws.Range('E47').Value = '  +15.69%  '
$ws.Range('D48').Value = '2.024.80'
$ws.Range('E48').Value = '  -4.11%  '
$ws.Range('D50').Value = '3.19'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('D51').Value = '0.0323'
$ws.Range('E51').Value = '  -2.34%  '
